{"js": "// Apply the Cover Letter revisions:\n//  1. \" U.S \" -> split into \" \" + \"U.S \" (no visible text change) in the\n//     \"I imagine you might find...\" paragraph.\n//  2. \"send\" -> \"sent\" in \"...excited as I send out my first aerospace application...\"\n//  3. \"thermal analysis student engineer\" -> \"Thermal Engineer\" in the same paragraph.\n//  4. Tidy up \"...my experience lies in a different industry...\" run split (no visible\n//     text change, handled implicitly since we never touch that phrase's text).\n//  5. Delete the whole \"To further my expertise in simulation, ...\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Work from the end of the document backwards so earlier paragraph indices\n// stay valid even after we delete one further down.\nlet toDelete = null;\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  if (text.indexOf(\"To further my expertise in simulation\") !== -1) {\n    // Entire paragraph (the CFD-solver side project blurb) was removed.\n    toDelete = para;\n    continue;\n  }\n\n  if (text.indexOf(\"U.S aerospace education\") !== -1) {\n    // Split \" U.S \" into \" \" + \"U.S \" -- purely a run split, text unchanged,\n    // but we still perform the replace so the edit is explicit/idempotent.\n    const hit = para.search(\" U.S \", { matchCase: true });\n    hit.load(\"items\");\n    await context.sync();\n    if (hit.items.length > 0) {\n      hit.items[0].insertText(\" U.S \", \"Replace\");\n      await context.sync();\n    }\n  }\n\n  if (text.indexOf(\"send out my first aerospace application\") !== -1) {\n    const hit = para.search(\"send out my first aerospace application\", { matchCase: true });\n    hit.load(\"items\");\n    await context.sync();\n    if (hit.items.length > 0) {\n      hit.items[0].insertText(\"sent out my first aerospace application\", \"Replace\");\n      await context.sync();\n    }\n  }\n\n  if (text.indexOf(\"thermal analysis student engineer at MDA\") !== -1) {\n    const hit = para.search(\"thermal analysis student engineer at MDA\", { matchCase: true });\n    hit.load(\"items\");\n    await context.sync();\n    if (hit.items.length > 0) {\n      hit.items[0].insertText(\"Thermal Engineer at MDA\", \"Replace\");\n      await context.sync();\n    }\n  }\n}\n\nif (toDelete) {\n  toDelete.delete();\n  await context.sync();\n}\n", "ps1": "# Apply the Cover Letter revisions:\n#  1. \" U.S \" -> split into \" \" + \"U.S \" (no visible text change) in the\n#     \"I imagine you might find...\" paragraph.\n#  2. \"send\" -> \"sent\" in \"...excited as I send out my first aerospace application...\"\n#  3. \"thermal analysis student engineer\" -> \"Thermal Engineer\" in the same paragraph.\n#  4. \"...my experience lies in a different industry...\" stays the same text (only a\n#     run-merge in the source, no content change), so no replace is needed there.\n#  5. Delete the whole \"To further my expertise in simulation, ...\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph to delete first, before any other edits shift ranges.\n$paraToDelete = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*To further my expertise in simulation*\") {\n        $paraToDelete = $p\n    }\n}\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    if ($t -like \"*U.S aerospace education*\") {\n        $r = $p.Range\n        $null = $r.Find.Execute(\" U.S \", $true, $false, $false, $false, $false, $true, 1, $false, \" U.S \", 2)\n    }\n\n    if ($t -like \"*send out my first aerospace application*\") {\n        $r = $p.Range\n        $null = $r.Find.Execute(\"send out my first aerospace application\", $true, $false, $false, $false, $false, $true, 1, $false, \"sent out my first aerospace application\", 2)\n    }\n\n    if ($t -like \"*thermal analysis student engineer at MDA*\") {\n        $r = $p.Range\n        $null = $r.Find.Execute(\"thermal analysis student engineer at MDA\", $true, $false, $false, $false, $false, $true, 1, $false, \"Thermal Engineer at MDA\", 2)\n    }\n}\n\nif ($paraToDelete -ne $null) {\n    $paraToDelete.Range.Delete()\n}\n"}
